$wb = $excel.ActiveWorkbook

# Existing sheets before the edit: [1]=2021-Q1, [2]=总计
$ws2021 = $wb.Worksheets.Item(1)
$wsTotal = $wb.Worksheets.Item(2)
$totalName = $wsTotal.Name

# --- Create the new "2022-Q1" sheet by duplicating "总计" so its header/row
# styles (s="2") are preserved exactly, then position it right after 2021-Q1.
# NOTE: inserting a sheet shifts index-based references, so the original
# "总计" handle ($wsTotal) must be re-fetched by name afterwards.
$wsTotal.Copy($null, $ws2021)
$ws2022 = $wb.Worksheets.Item(2)
$ws2022.Name = "2022-Q1"

# Re-fetch the original "总计" sheet by name now that indices have shifted.
$wsTotal = $wb.Worksheets.Item($totalName)

# Extend header styling (copy the existing D1 "s=2" formatting across E1:H1)
$ws2022.Range("D1").Copy()
$ws2022.Range("E1:H1").PasteSpecial(-4122)

# Extend the row index column style (copy A2 "s=2" formatting down to A3)
$ws2022.Range("A2").Copy()
$ws2022.Range("A3").PasteSpecial(-4122)

# --- Header row ---
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"

# --- Row 2: fund 007257 ---
$ws2022.Range("A2").Value = 0

$ws2022.Range("B2").NumberFormat = "@"
$ws2022.Range("B2").Value = "007257"
$ws2022.Range("B2").ClearFormats()

$ws2022.Range("C2").NumberFormat = "@"
$ws2022.Range("C2").Value = "凯石沣混合A"
$ws2022.Range("C2").ClearFormats()

$ws2022.Range("D2").NumberFormat = "@"
$ws2022.Range("D2").Value = "0.17"
$ws2022.Range("D2").ClearFormats()

$ws2022.Range("E2").NumberFormat = "@"
$ws2022.Range("E2").Value = "72.94"
$ws2022.Range("E2").ClearFormats()

$ws2022.Range("F2").NumberFormat = "@"
$ws2022.Range("F2").Value = "3.80"
$ws2022.Range("F2").ClearFormats()

$ws2022.Range("G2").NumberFormat = "@"
$ws2022.Range("G2").Value = "0.0065"
$ws2022.Range("G2").ClearFormats()

$ws2022.Range("H2").Value = 2

# --- Row 3: fund 007258 ---
$ws2022.Range("A3").Value = 1

$ws2022.Range("B3").NumberFormat = "@"
$ws2022.Range("B3").Value = "007258"
$ws2022.Range("B3").ClearFormats()

$ws2022.Range("C3").NumberFormat = "@"
$ws2022.Range("C3").Value = "凯石沣混合C"
$ws2022.Range("C3").ClearFormats()

$ws2022.Range("D3").NumberFormat = "@"
$ws2022.Range("D3").Value = "0.05"
$ws2022.Range("D3").ClearFormats()

$ws2022.Range("E3").NumberFormat = "@"
$ws2022.Range("E3").Value = "72.94"
$ws2022.Range("E3").ClearFormats()

$ws2022.Range("F3").NumberFormat = "@"
$ws2022.Range("F3").Value = "3.80"
$ws2022.Range("F3").ClearFormats()

$ws2022.Range("G3").NumberFormat = "@"
$ws2022.Range("G3").Value = "0.0019"
$ws2022.Range("G3").ClearFormats()

$ws2022.Range("H3").Value = 2

# --- Update the "总计" sheet: push the existing 2021-Q1 row down to row 3,
# then write the new 2022-Q1 summary row into row 2. Format + values are
# copied in separate passes; a plain xlPasteAll does not reliably carry the
# row-index column's style (s="2") over to the new row in this host. ---
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4122)  # xlPasteFormats
$wsTotal.Range("A2:D2").Copy()
$wsTotal.Range("A3:D3").PasteSpecial(-4163)  # xlPasteValues
$wsTotal.Range("A3").Value = 1

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01
